# Update the table style used by the three data tables (slides 14, 15, 16)
# from "{8DFB46A2-3820-4D54-A9E1-2525D23C5124}" to
# "{ABA3DB55-FA22-48EC-A4E7-299BF17D6632}".
#
# In each of these slides the table (a p:graphicFrame) is the first shape
# in the shape tree, so Shapes.Item(1) is the table shape.

$p = $ppt.ActivePresentation

$newStyleId = "{ABA3DB55-FA22-48EC-A4E7-299BF17D6632}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    $tableShape = $slide.Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle($newStyleId, $false)
    }
}
